$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: UserID values (text), Column B: Password values (numbers)
$userIds = @("a1", "a2", "a3", "a4", "a5", "a6")
$passwords = @(11, 22, 33, 44, 55, 66)

for ($i = 0; $i -lt $userIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $userIds[$i]
    $ws.Cells.Item($row, 2).Value = $passwords[$i]
}

$ws.Range("C3").Select()
